# IOT_dummy.xlsx edit: add a new "units" worksheet describing the unit of
# measure for each row category, positioned right before the "usefull"
# sheet, and nudge a couple of unrelated selection/active-cell bits that
# moved around as a side effect of the author's interactive session.

$wb = $excel.ActiveWorkbook

# --- selection tweaks on two pre-existing sheets (no sheet activation) ---
$wsAllData = $wb.Worksheets.Item("all_data")
$wsAllData.Range("L19").Select()

$wsE = $wb.Worksheets.Item("E")
$wsE.Range("E16").Select()

# --- insert the new "units" sheet right before "usefull" ---
$wsUsefull = $wb.Worksheets.Item("usefull")
$units = $wb.Worksheets.Add($wsUsefull)
$units.Name = "units"

# Fill the "unit" column first (top to bottom), matching the order the
# shared-string table picked up new entries in the authored workbook.
$units.Range("C1").Value = "unit"
$units.Range("C2").Value = "eur"
$units.Range("C3").Value = "eur"
$units.Range("C4").Value = "eur"
$units.Range("C5").Value = "eur"
$units.Range("C6").Value = "eur"
$units.Range("C7").Value = "eur"
$units.Range("C8").Value = "kg"
$units.Range("C9").Value = "m2"

# Then the category column (A)
$units.Range("A2").Value = "Sector"
$units.Range("A3").Value = "Sector"
$units.Range("A4").Value = "Sector"
$units.Range("A5").Value = "Sector"
$units.Range("A6").Value = "Factor of production"
$units.Range("A7").Value = "Factor of production"
$units.Range("A8").Value = "Satellite account"
$units.Range("A9").Value = "Satellite account"

# Then the row-label column (B)
$units.Range("B2").Value = "sec1"
$units.Range("B3").Value = "sec2"
$units.Range("B4").Value = "sec1"
$units.Range("B5").Value = "sec2"
$units.Range("B6").Value = "Labor"
$units.Range("B7").Value = "Capital"
$units.Range("B8").Value = "Emission"
$units.Range("B9").Value = "Land"

# make it the active sheet / tab, with the final interactive selection
$units.Activate()
$units.Range("G8").Select()
